$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp shown in the title row (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Septiembre de 2020 a las 06:29"

# India (row 5): Casos totales, Nuevos casos, Casos activos, Recuperados
$ws.Range("B5").Value = 5818570
$ws.Range("C5").Value = 2467
$ws.Range("D5").Value = 4756164
$ws.Range("E5").Value = 970089

# Tailandia (row 139): Casos totales, Nuevos casos, Casos activos, Recuperados
$ws.Range("B139").Value = 3519
$ws.Range("C139").Value = 3
$ws.Range("D139").Value = 3360
$ws.Range("E139").Value = 100

# Mongolia (row 186): Casos activos, Recuperados
$ws.Range("D186").Value = 303
$ws.Range("E186").Value = 10

# Butan (row 188): Casos totales, Nuevos casos, Casos activos
$ws.Range("B188").Value = 263
$ws.Range("C188").Value = 2
$ws.Range("D188").Value = 198

# Rows 215 & 216 swap places in the updated country ranking (Islas Malvinas
# now ranks above Montserrat), so their full data rows are exchanged.
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 1
